{"js": "// Remove the old instructional note block (\"HO\u00c0N CH\u1ec8NH THI\u1ebeT K\u1ebe D\u1eee LI\u1ec6U\n// T\u1ea0I \u0110\u00c2Y.\", \"=> CH\u00da \u00dd:\", and the \"- D\u1ef0A TR\u00caN ... / - \u0110\u1ea2M B\u1ea2O ...\" notes),\n// leaving only the trailing blank paragraph before the section break \u2014\n// per commit \"Ho\u00e0n th\u00e0nh thi\u1ebft k\u1ebf d\u1eef li\u1ec7u.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify the contiguous run of paragraphs that make up the note block by\n// matching their text against the known content being removed, then delete\n// exactly that run (keeps any other paragraph \u2014 e.g. the trailing blank\n// line that stays \u2014 untouched).\nconst markers = [\n  \"HO\u00c0N CH\u1ec8NH THI\u1ebeT K\u1ebe D\u1eee LI\u1ec6U T\u1ea0I \u0110\u00c2Y.\",\n  \"=> CH\u00da \u00dd:\",\n  \"- D\u1ef0A TR\u00caN Y\u00caU C\u1ea6U TH\u1ea6Y \u0110\u01afA RA V\u00c0O NG\u00c0Y TH\u1ee8 2 (17/05).\",\n];\n\nconst items = paragraphs.items;\nlet firstIndex = -1;\nlet lastIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (markers.some((marker) => text.startsWith(marker))) {\n    if (firstIndex === -1) firstIndex = i;\n    lastIndex = i;\n  }\n}\n\nif (firstIndex !== -1) {\n  // Collect the whole matched run (this also folds in the blank separator\n  // paragraph that sits right after the opening note) and delete it as one\n  // unit, end first so earlier indices stay valid while the loop runs.\n  const toDelete = [];\n  for (let i = firstIndex; i <= lastIndex; i++) {\n    toDelete.push(items[i]);\n  }\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the old instructional note block (\"HO\u00c0N CH\u1ec8NH THI\u1ebeT K\u1ebe D\u1eee LI\u1ec6U\n# T\u1ea0I \u0110\u00c2Y.\", \"=> CH\u00da \u00dd:\", and the \"- D\u1ef0A TR\u00caN ... / - \u0110\u1ea2M B\u1ea2O ...\" notes),\n# leaving only the trailing blank paragraph before the section break \u2014\n# per commit \"Ho\u00e0n th\u00e0nh thi\u1ebft k\u1ebf d\u1eef li\u1ec7u.\"\n$d = $word.ActiveDocument\n\n$markers = @(\n    \"HO\u00c0N CH\u1ec8NH THI\u1ebeT K\u1ebe D\u1eee LI\u1ec6U T\u1ea0I \u0110\u00c2Y.\",\n    \"=> CH\u00da \u00dd:\",\n    \"- D\u1ef0A TR\u00caN Y\u00caU C\u1ea6U TH\u1ea6Y \u0110\u01afA RA V\u00c0O NG\u00c0Y TH\u1ee8 2 (17/05).\"\n)\n\n# Find the contiguous run of paragraphs whose text matches the note block\n# so the whole block is removed as a unit, wherever it sits in the body.\n$firstIndex = -1\n$lastIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]11, [char]7)\n    foreach ($marker in $markers) {\n        if ($t.StartsWith($marker)) {\n            if ($firstIndex -eq -1) { $firstIndex = $i }\n            $lastIndex = $i\n        }\n    }\n}\n\nif ($firstIndex -ne -1) {\n    $rangeStart = $d.Paragraphs.Item($firstIndex).Range.Start\n    $rangeEnd = $d.Paragraphs.Item($lastIndex).Range.End\n    $deleteRange = $d.Range($rangeStart, $rangeEnd)\n    $deleteRange.Delete()\n}\n"}
